$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $null = $d.Content.Find.Execute($find, $false, $false, $false, $false, $false, $true, 1, $false, $replace, 2)
}

# "Simple Tasks" section
Replace-Text "Issue (1):"  "Issue (1) -:"
Replace-Text "Issue (2):"  "Issue (2) _:"
Replace-Text "Issue (4):"  "Issue (4) _:"
Replace-Text "Issue (5):"  "Issue (5) -:"
Replace-Text "Issue (6):"  "Issue (6) -:"
Replace-Text "Issue (11):" "Issue (11) -:"
Replace-Text "Issue (15):" "Issue (15) -:"

# "Moderate Tasks" section
Replace-Text "Issue (3):"  "Issue (3) _:"
Replace-Text "Issue (8):"  "Issue (8) -:"
Replace-Text "Issue (10):" "Issue (10) #:"
Replace-Text "Issue (12):" "Issue (12) -:"
Replace-Text "Issue (13):" "Issue (13) #:"
Replace-Text "Issue (14):" "Issue (14) -:"

# "Tasks I didn't fully understand" section
Replace-Text "Issue (19) & Issue (20) & Issue (21):" "Issue (13) & Issue (19) & Issue (20) & Issue (21) #:"
